$d = $word.ActiveDocument

# --- Update the date shown in the first-page header (2023-09-13 -> 2023-09-15) ---
$sec = $d.Sections.Item(1)
$hdrs = $sec.Headers
for ($hi = 1; $hi -le $hdrs.Count; $hi++) {
    $hdr = $hdrs.Item($hi)
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
    }
}

# --- Append new 'Knärot' section + references after 'BILAGA 1 - Fridlysta arter' ---
# Phase 1: create all the empty paragraphs first (so that later per-paragraph
# formatting operations on one paragraph cannot bleed into a not-yet-created one).
$anchor = $d.Paragraphs.Last
$anchor.Range.InsertParagraphAfter()
$p0 = $d.Paragraphs.Last
$anchor = $p0
$anchor.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$anchor = $p1
$anchor.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$anchor = $p2
$anchor.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$anchor = $p3
$anchor.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$anchor = $p4
$anchor.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Last
$anchor = $p5
$anchor.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Last
$anchor = $p6
$anchor.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Last
$anchor = $p7
$anchor.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs.Last
$anchor = $p8
$anchor.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs.Last
$anchor = $p9
$anchor.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs.Last
$anchor = $p10
$anchor.Range.InsertParagraphAfter()
$p11 = $d.Paragraphs.Last
$anchor = $p11
$anchor.Range.InsertParagraphAfter()
$p12 = $d.Paragraphs.Last
$anchor = $p12

# Phase 2: set the style + full plain text of each paragraph.
$p0.Style = "Heading1"
$p0.Range.Text = 'Knärot – ekologi samt krav på livsmiljön'
$p1.Style = "Normal"
$p1.Range.Text = 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).'
$p2.Style = "Normal"
$p2.Range.Text = 'Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”'
$p3.Style = "Normal"
$p3.Range.Text = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”'
$p4.Style = "Normal"
$p4.Range.Text = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).'
$p5.Style = "Normal"
$p5.Range.Text = 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).'
$p6.Style = "Heading2"
$p6.Range.Text = 'Referenser - knärot'
$p7.Style = "Normal"
$p7.Range.Text = 'de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025'
$p8.Style = "Normal"
$p8.Range.Text = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 '
$p9.Style = "Normal"
$p9.Range.Text = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853'
$p10.Style = "Normal"
$p10.Range.Text = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.'
$p11.Style = "Normal"
$p11.Range.Text = 'Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/'
$p12.Style = "Normal"
$p12.Range.Text = 'SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala '

# Phase 3: apply italic formatting to the relevant sub-ranges of each paragraph.
$s2 = $p2.Range.Start
$d.Range($s2 + 34, $s2 + 116).Font.Italic = $true
$d.Range($s2 + 278, $s2 + 483).Font.Italic = $true
$d.Range($s2 + 490, $s2 + 608).Font.Italic = $true
$s3 = $p3.Range.Start
$d.Range($s3 + 205, $s3 + 1070).Font.Italic = $true
$s7 = $p7.Range.Start
$d.Range($s7 + 33, $s7 + 113).Font.Italic = $true
$s8 = $p8.Range.Start
$d.Range($s8 + 62, $s8 + 176).Font.Italic = $true
$s9 = $p9.Range.Start
$d.Range($s9 + 117, $s9 + 207).Font.Italic = $true
$s10 = $p10.Range.Start
$d.Range($s10 + 54, $s10 + 121).Font.Italic = $true
$s11 = $p11.Range.Start
$d.Range($s11 + 22, $s11 + 57).Font.Italic = $true
$s12 = $p12.Range.Start
$d.Range($s12 + 25, $s12 + 61).Font.Italic = $true

